$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.7028603333333333
$ws.Range("H2").Value = 2.108581
$ws.Range("I2").Value = 0.003742168187771943
$ws.Range("J2").Value = 0.003742168187771943
$ws.Range("M2").Value = 71.05094633333333
$ws.Range("N2").Value = 213.152839
$ws.Range("O2").Value = 0.8240565632932695
$ws.Range("P2").Value = 0.8240565632932696
$ws.Range("Q2").Value = 49.93889182349544
$ws.Range("R2").Value = 449.450026411459
$ws.Range("S2").Value = 0.003083758256080749
$ws.Range("T2").Value = 0.00308375825608075
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.7028603333333333
$ws.Range("H3").Value = 2.108581
$ws.Range("I3").Value = 0.003742168187771943
$ws.Range("J3").Value = 0.003742168187771943
$ws.Range("O3").Value = 0.1323102827659759
$ws.Range("P3").Value = 0.132310282765976
$ws.Range("Q3").Value = 8.018173985266445
$ws.Range("R3").Value = 72.163565867398
$ws.Range("S3").Value = 0.0004951273310819454
$ws.Range("T3").Value = 0.0004951273310819455
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.7028603333333333
$ws.Range("H4").Value = 2.108581
$ws.Range("I4").Value = 0.003742168187771943
$ws.Range("J4").Value = 0.003742168187771943
$ws.Range("M4").Value = 3.762092333333333
$ws.Range("N4").Value = 11.286277
$ws.Range("O4").Value = 0.04363315394075456
$ws.Range("P4").Value = 0.04363315394075455
$ws.Range("Q4").Value = 2.644225471437444
$ws.Range("R4").Value = 23.798029242937
$ws.Range("S4").Value = 0.0001632826006092477
$ws.Range("T4").Value = 0.0001632826006092477
$ws.Range("I5").Value = 0.9809344360866079
$ws.Range("J5").Value = 0.9809344360866079
$ws.Range("M5").Value = 71.05094633333333
$ws.Range("N5").Value = 213.152839
$ws.Range("O5").Value = 0.8240565632932695
$ws.Range("P5").Value = 0.8240565632932696
$ws.Range("Q5").Value = 13090.4802327543
$ws.Range("R5").Value = 117814.3220947887
$ws.Range("S5").Value = 0.8083454602175514
$ws.Range("T5").Value = 0.8083454602175515
$ws.Range("I6").Value = 0.9809344360866079
$ws.Range("J6").Value = 0.9809344360866079
$ws.Range("O6").Value = 0.1323102827659759
$ws.Range("P6").Value = 0.132310282765976
$ws.Range("S6").Value = 0.1297877126135022
$ws.Range("T6").Value = 0.1297877126135022
$ws.Range("I7").Value = 0.9809344360866079
$ws.Range("J7").Value = 0.9809344360866079
$ws.Range("M7").Value = 3.762092333333333
$ws.Range("N7").Value = 11.286277
$ws.Range("O7").Value = 0.04363315394075456
$ws.Range("P7").Value = 0.04363315394075455
$ws.Range("Q7").Value = 693.1307444133528
$ws.Range("R7").Value = 6238.176699720175
$ws.Range("S7").Value = 0.04280126325555423
$ws.Range("T7").Value = 0.04280126325555422
$ws.Range("G8").Value = 2.258412
$ws.Range("H8").Value = 6.775236
$ws.Range("I8").Value = 0.01202423460319866
$ws.Range("J8").Value = 0.01202423460319866
$ws.Range("M8").Value = 71.05094633333333
$ws.Range("N8").Value = 213.152839
$ws.Range("O8").Value = 0.8240565632932695
$ws.Range("P8").Value = 0.8240565632932696
$ws.Range("Q8").Value = 160.462309810556
$ws.Range("R8").Value = 1444.160788295004
$ws.Range("S8").Value = 0.009908649443343895
$ws.Range("T8").Value = 0.009908649443343895
$ws.Range("G9").Value = 2.258412
$ws.Range("H9").Value = 6.775236
$ws.Range("I9").Value = 0.01202423460319866
$ws.Range("J9").Value = 0.01202423460319866
$ws.Range("O9").Value = 0.1323102827659759
$ws.Range("P9").Value = 0.132310282765976
$ws.Range("Q9").Value = 25.763781917432
$ws.Range("R9").Value = 231.874037256888
$ws.Range("S9").Value = 0.001590929880393646
$ws.Range("T9").Value = 0.001590929880393647
$ws.Range("G10").Value = 2.258412
$ws.Range("H10").Value = 6.775236
$ws.Range("I10").Value = 0.01202423460319866
$ws.Range("J10").Value = 0.01202423460319866
$ws.Range("M10").Value = 3.762092333333333
$ws.Range("N10").Value = 11.286277
$ws.Range("O10").Value = 0.04363315394075456
$ws.Range("P10").Value = 0.04363315394075455
$ws.Range("Q10").Value = 8.496354470707999
$ws.Range("R10").Value = 76.46719023637199
$ws.Range("S10").Value = 0.0005246552794611148
$ws.Range("T10").Value = 0.0005246552794611147
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.619654
$ws.Range("H11").Value = 1.858962
$ws.Range("I11").Value = 0.003299161122421622
$ws.Range("J11").Value = 0.003299161122421622
$ws.Range("M11").Value = 71.05094633333333
$ws.Range("N11").Value = 213.152839
$ws.Range("O11").Value = 0.8240565632932695
$ws.Range("P11").Value = 0.8240565632932696
$ws.Range("Q11").Value = 44.02700309923533
$ws.Range("R11").Value = 396.243027893118
$ws.Range("S11").Value = 0.002718695376293528
$ws.Range("T11").Value = 0.002718695376293528
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.619654
$ws.Range("H12").Value = 1.858962
$ws.Range("I12").Value = 0.003299161122421622
$ws.Range("J12").Value = 0.003299161122421622
$ws.Range("O12").Value = 0.1323102827659759
$ws.Range("P12").Value = 0.132310282765976
$ws.Range("Q12").Value = 7.068962846577334
$ws.Range("R12").Value = 63.620665619196
$ws.Range("S12").Value = 0.0004365129409981194
$ws.Range("T12").Value = 0.0004365129409981194
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.619654
$ws.Range("H13").Value = 1.858962
$ws.Range("I13").Value = 0.003299161122421622
$ws.Range("J13").Value = 0.003299161122421622
$ws.Range("M13").Value = 3.762092333333333
$ws.Range("N13").Value = 11.286277
$ws.Range("O13").Value = 0.04363315394075456
$ws.Range("P13").Value = 0.04363315394075455
$ws.Range("Q13").Value = 2.331195562719333
$ws.Range("R13").Value = 20.980760064474
$ws.Range("S13").Value = 0.0001439528051299753
$ws.Range("T13").Value = 0.0001439528051299752
